$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values (number of reps / samples per subject)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 (CON) passive torque values
$ws.Range("B2").Value = 17.475217118104748
$ws.Range("C2").Value = 46.915574224234724
$ws.Range("D2").Value = 17.815848497584529
$ws.Range("E2").Value = 33.959880252933111

# Update row 3 (STR) passive torque values
$ws.Range("B3").Value = 21.561630270302388
$ws.Range("C3").Value = 25.305169314230415
$ws.Range("D3").Value = 17.426995298519632
$ws.Range("E3").Value = 21.964335007884028

# Update selection to match new narrower range of interest
$ws.Range("B1:E3").Select()
